$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new portfolio data row (row 28) following the same layout as
# the existing rows: col A is a literal text date, B/C/D are numbers.
$row = 28

# Force column A to be stored as text (not auto-converted to a date serial),
# then clear the format so the cell ends up with no explicit style - matching
# the other plain data rows above it.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-12"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 57.11000061035156
$ws.Cells.Item($row, 3).Value = 715.25
$ws.Cells.Item($row, 4).Value = 321.3999938964844
